$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New headers for columns I (I0) and J (IF), matching the style of the
# existing header row (bold font, thin border, centered alignment).
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

$headerRange = $ws.Range("I1:J1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160
$headerRange.Borders.LineStyle = 1
$headerRange.Borders.Weight = 2

# Row, I0, IF
$data = @(
    @(2,8,8),
    @(3,9,9),
    @(4,7,7),
    @(5,8,8),
    @(6,9,9),
    @(7,8,8),
    @(8,9,9),
    @(9,9,9),
    @(10,9,9),
    @(11,9,10),
    @(12,9,9),
    @(13,9,9),
    @(14,9,9),
    @(15,9,9),
    @(16,8,8),
    @(17,9,9),
    @(18,9,9),
    @(19,8,8),
    @(20,8,8),
    @(21,7,8),
    @(22,9,9),
    @(23,8,8),
    @(24,8,8),
    @(25,8,8),
    @(26,8,8),
    @(27,9,9),
    @(28,8,8),
    @(29,8,8),
    @(30,9,9),
    @(31,8,8),
    @(32,9,9),
    @(33,8,8),
    @(34,7,7),
    @(35,7,8),
    @(36,8,8),
    @(37,8,8),
    @(38,8,8),
    @(39,8,8),
    @(40,8,8),
    @(41,7,8),
    @(42,9,9),
    @(43,7,8),
    @(44,8,8),
    @(45,7,8),
    @(46,7,8),
    @(47,7,8),
    @(48,8,8),
    @(49,8,8),
    @(50,7,8),
    @(51,8,8),
    @(52,8,8),
    @(53,7,8),
    @(54,7,8),
    @(55,9,9),
    @(56,7,8),
    @(57,8,8),
    @(58,8,8),
    @(59,8,8),
    @(60,9,9),
    @(61,8,8),
    @(62,9,9),
    @(63,7,8),
    @(64,8,9),
    @(65,7,8),
    @(66,9,9),
    @(67,8,9),
    @(68,8,8),
    @(69,7,8),
    @(70,9,9),
    @(71,7,8),
    @(72,8,8),
    @(73,8,8),
    @(74,9,9),
    @(75,7,8),
    @(76,9,9),
    @(77,7,8),
    @(78,4,4),
    @(79,9,9),
    @(80,8,8),
    @(81,7,7),
    @(82,4,4)
)

foreach ($row in $data) {
    $r = $row[0]
    $i0 = $row[1]
    $iff = $row[2]
    $ws.Cells.Item($r, 9).Value = $i0
    $ws.Cells.Item($r, 10).Value = $iff
}
